# Applies the 2023-12-26 cryptos-list price/volume refresh described by
# the commit diff: most rows get updated Price (D) / Volume(1h) (E) text,
# and row 51 swaps from "Cronos" to "Aave" (new name/link/price/volume).
#
# The D/E columns are stored as literal text (t="inlineStr") in the
# workbook, not numbers - e.g. "42.688.65" or "  -1.47%  ". Assigning a
# plain numeric-looking string to Range.Value lets Excel's COM layer
# auto-coerce it to a real number, and flipping NumberFormat to force text
# normally mints a brand-new cell style. Set-TextValue below avoids both:
# it marks the cell as text before the write (so the value is kept
# verbatim as a string) and then restores the cell's original Style
# afterwards so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$edits = @(

    @{ Ref = "D2"; Value = '42.688.65' },
    @{ Ref = "E2"; Value = '  -1.47%  ' },
    @{ Ref = "D3"; Value = '2.242.06' },
    @{ Ref = "E3"; Value = '  -1.67%  ' },
    @{ Ref = "D5"; Value = '115.33' },
    @{ Ref = "E5"; Value = '  +2.85%  ' },
    @{ Ref = "D6"; Value = '288.66' },
    @{ Ref = "E6"; Value = '  +9.05%  ' },
    @{ Ref = "D7"; Value = '0.629' },
    @{ Ref = "E7"; Value = '  -2.91%  ' },
    @{ Ref = "E8"; Value = '  -0.03%  ' },
    @{ Ref = "D9"; Value = '0.617' },
    @{ Ref = "E9"; Value = '  +1.16%  ' },
    @{ Ref = "E10"; Value = '  +0.31%  ' },
    @{ Ref = "D11"; Value = '0.0934' },
    @{ Ref = "E11"; Value = '  -0.20%  ' },
    @{ Ref = "D12"; Value = '9.22' },
    @{ Ref = "E12"; Value = '  +0.16%  ' },
    @{ Ref = "E13"; Value = '  -2.86%  ' },
    @{ Ref = "D14"; Value = '15.45' },
    @{ Ref = "E14"; Value = '  +0.99%  ' },
    @{ Ref = "D15"; Value = '0.890' },
    @{ Ref = "E15"; Value = '  +3.37%  ' },
    @{ Ref = "D16"; Value = '2.580.30' },
    @{ Ref = "E16"; Value = '  -1.66%  ' },
    @{ Ref = "D17"; Value = '2.249.73' },
    @{ Ref = "E17"; Value = '  -1.22%  ' },
    @{ Ref = "D18"; Value = '42.754.26' },
    @{ Ref = "E18"; Value = '  -1.31%  ' },
    @{ Ref = "E19"; Value = '  -1.31%  ' },
    @{ Ref = "D20"; Value = '7.24' },
    @{ Ref = "E20"; Value = '  +7.81%  ' },
    @{ Ref = "D21"; Value = '73.46' },
    @{ Ref = "E21"; Value = '  +1.54%  ' },
    @{ Ref = "D22"; Value = '3.34' },
    @{ Ref = "E22"; Value = '  +16.12%  ' },
    @{ Ref = "D23"; Value = '2.38' },
    @{ Ref = "E23"; Value = '  -1.96%  ' },
    @{ Ref = "D24"; Value = '232.44' },
    @{ Ref = "E24"; Value = '  -0.93%  ' },
    @{ Ref = "D25"; Value = '9.22' },
    @{ Ref = "E25"; Value = '  -2.28%  ' },
    @{ Ref = "D26"; Value = '12.17' },
    @{ Ref = "E26"; Value = '  +5.20%  ' },
    @{ Ref = "E27"; Value = '  -1.42%  ' },
    @{ Ref = "E28"; Value = '  -0.66%  ' },
    @{ Ref = "E29"; Value = '  -1.58%  ' },
    @{ Ref = "E30"; Value = '  -1.41%  ' },
    @{ Ref = "D31"; Value = '2.19' },
    @{ Ref = "E31"; Value = '  -2.28%  ' },
    @{ Ref = "D32"; Value = '175.32' },
    @{ Ref = "E32"; Value = '  +1.09%  ' },
    @{ Ref = "D33"; Value = '21.29' },
    @{ Ref = "E33"; Value = '  -1.40%  ' },
    @{ Ref = "D34"; Value = '0.0912' },
    @{ Ref = "E34"; Value = '  +1.11%  ' },
    @{ Ref = "D35"; Value = '4.54' },
    @{ Ref = "E35"; Value = '  +18.57%  ' },
    @{ Ref = "D36"; Value = '5.60' },
    @{ Ref = "E36"; Value = '  -0.57%  ' },
    @{ Ref = "E37"; Value = '  -2.41%  ' },
    @{ Ref = "E38"; Value = '  -1.45%  ' },
    @{ Ref = "D39"; Value = '4.66' },
    @{ Ref = "E39"; Value = '  -0.96%  ' },
    @{ Ref = "D40"; Value = '0.106' },
    @{ Ref = "E40"; Value = '  +2.38%  ' },
    @{ Ref = "E41"; Value = '  +3.01%  ' },
    @{ Ref = "D42"; Value = '73.04' },
    @{ Ref = "E42"; Value = '  -1.98%  ' },
    @{ Ref = "D43"; Value = '13.52' },
    @{ Ref = "E43"; Value = '  -6.47%  ' },
    @{ Ref = "D44"; Value = '0.237' },
    @{ Ref = "E44"; Value = '  +0.35%  ' },
    @{ Ref = "E45"; Value = '  +0.12%  ' },
    @{ Ref = "E46"; Value = '  -1.09%  ' },
    @{ Ref = "D47"; Value = '5.62' },
    @{ Ref = "E47"; Value = '  -6.77%  ' },
    @{ Ref = "D48"; Value = '1.32' },
    @{ Ref = "E48"; Value = '  +4.11%  ' },
    @{ Ref = "D49"; Value = '8.55' },
    @{ Ref = "E49"; Value = '  -0.20%  ' },
    @{ Ref = "D50"; Value = '0.653' },
    @{ Ref = "E50"; Value = '  +6.02%  ' },
    @{ Ref = "B51"; Value = 'Aave' },
    @{ Ref = "C51"; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Ref = "D51"; Value = '101.36' },
    @{ Ref = "E51"; Value = '  +1.13%  ' }

)

foreach ($e in $edits) {
    Set-TextValue $ws $e.Ref $e.Value
}
